# Adds a new module entry ("天外流星" by 玉米man) as row 48 on Sheet1.
# Sheet2 already carries formulas through row 52 that reference Sheet1 rows,
# so they will recompute automatically once Sheet1!A48:Z48 is populated.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---- Sheet1 row 48 ---------------------------------------------------
$ws1.Range("A48").Value = '天外流星'
$ws1.Range("B48").Value = '玉米man'
$ws1.Range("C48").Value = 'DND5E'
$ws1.Range("D48").Value = '短模组（开阔世界）'
$ws1.Range("E48").Value = '短篇'
$ws1.Range("F48").Value = 1
$ws1.Range("G48").Value = 45537
$ws1.Range("H48").Value = '不定'
$ws1.Range("I48").Value = 4
$ws1.Range("J48").Value = 4
$ws1.Range("K48").Value = 'T2'
$ws1.Range("L48").Value = 7
$ws1.Range("M48").Value = 7
$ws1.Range("N48").Value = 7
$ws1.Range("O48").Value = 9
$ws1.Range("P48").Value = '酒馆中寻找姐姐的小女孩，在老板的帮助下发布委托。冒险者接下委托，并前往地城直面恐怖的真相。'
$ws1.Range("Q48").Value = '第54期逸闻酒馆活动'
$ws1.Range("R48").Value = '无'
$ws1.Range("S48").Value = '【彗星】【虚假记忆】【喧闹】【纸牌】'
$ws1.Range("T48").Value = '有'
$ws1.Range("U48").Value = '有'
$ws1.Range("V48").Value = '有'
$ws1.Range("W48").Value = '无'
$ws1.Range("X48").Value = '冒险者们在酒馆歇息，一个小女孩进入了酒馆，她四处环顾，怯生生的走到了吧台。正擦着酒杯的吧台老板虽显疑惑，但仍为小女孩递上一杯橙汁，告诉小女孩这里可不是孩子该来的地方。小女孩有些害怕，在喝了几口橙汁之后用含糊不清的话小声说到："我要……委托"'
$ws1.Range("Z48").Value = 'D&D5E, T2(5-10), 侦探故事, 冒险故事, 短篇, 第54期 喧闹的纸牌彗星带来虚假记忆'

# Y48 carries a hyperlink out to the module's write-up; adding the
# hyperlink also fills the visible cell text with the URL.
$ws1.Hyperlinks.Add($ws1.Range("Y48"), 'https://wd.52hertz.work/?p=609') | Out-Null
# Match the hyperlink style already used by the other Y-column links
# (Hyperlinks.Add leaves behind a near-duplicate style otherwise).
$ws1.Range("Y47").Copy() | Out-Null
$ws1.Range("Y48").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 48 wraps onto multiple lines like the rows above it.
$ws1.Rows.Item(48).RowHeight = 54

# ---- view state (cosmetic, best effort) -------------------------------
$ws2.Activate()
$ws2.Range("B47").Select()

$ws1.Activate()
$ws1.Range("X52").Select()

Write-Host "Row 48 added to Sheet1"
